$p = $ppt.ActivePresentation

# EMU -> points conversion (914400 EMU per inch, 72 points per inch)
function EMU-ToPoints($emu) {
    return $emu / 914400.0 * 72.0
}

# 1. Update the datetimeFigureOut placeholder text from 4/14/2021 to 4/15/2021
#    on the Slide Master and on every Slide Layout.
$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $sh = $sm.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "4/14/2021") {
            $sh.TextFrame.TextRange.Text = "4/15/2021"
        }
    }
}

for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $cl = $sm.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "4/14/2021") {
                $sh.TextFrame.TextRange.Text = "4/15/2021"
            }
        }
    }
}

# 2. On slide 4, rename the "owner" textbox to "performer" and move/resize it
#    (v3 ProcedureRequest "owner" -> v4 ServiceRequest "performer").
$s4 = $p.Slides.Item(4)
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $sh = $s4.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "owner") {
            $sh.TextFrame.TextRange.Text = "performer"
            $sh.Left = EMU-ToPoints(4588610)
            $sh.Top = EMU-ToPoints(5396346)
            $sh.Width = EMU-ToPoints(789487)
            $sh.Height = EMU-ToPoints(261610)
        }
    }
}
